$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04935063238890791
$ws.Range("D2").Value = 0.02329706158598555
$ws.Range("E2").Value = 0.09572615734814605
$ws.Range("F2").Value = 0.6575948666205704
$ws.Range("G2").Value = 0.002414832999980501
$ws.Range("K2").Value = 1.190993477032976
$ws.Range("M2").Value = 0.3848367263253465
$ws.Range("O2").Value = 2.220079531345846
$ws.Range("C3").Value = 0.04381608568382944
$ws.Range("D3").Value = 0.02184499405163365
$ws.Range("E3").Value = 0.09175192719149194
$ws.Range("F3").Value = 0.6570273404666551
$ws.Range("G3").Value = 0.002418046712899343
$ws.Range("K3").Value = 1.043158600739957
$ws.Range("M3").Value = 0.3435600473777072
$ws.Range("O3").Value = 2.234118524495727
$ws.Range("C4").Value = 0.04043264710733752
$ws.Range("D4").Value = 0.02095091101128332
$ws.Range("E4").Value = 0.08942706720871385
$ws.Range("F4").Value = 0.6573115144528856
$ws.Range("G4").Value = 0.002420123526819381
$ws.Range("K4").Value = 0.9521415680187033
$ws.Range("M4").Value = 0.3182726918021928
$ws.Range("O4").Value = 2.244907689826732
$ws.Range("C5").Value = 0.03905755951444689
$ws.Range("D5").Value = 0.02058596368797794
$ws.Range("E5").Value = 0.08850834500567473
$ws.Range("F5").Value = 0.6575858425592784
$ws.Range("G5").Value = 0.002420995973066664
$ws.Range("K5").Value = 0.9149914020390497
$ws.Range("M5").Value = 0.3079821650989487
$ws.Range("O5").Value = 2.249847838789037
$ws.Range("C6").Value = 0.03882944935303101
$ws.Range("D6").Value = 0.02052532902470006
$ws.Range("E6").Value = 0.08835751443313811
$ws.Range("F6").Value = 0.657640952726517
$ws.Range("G6").Value = 0.002421142422478213
$ws.Range("K6").Value = 0.9088190767143942
$ws.Range("M6").Value = 0.3062742940158429
$ws.Range("O6").Value = 2.25070091343585
$ws.Range("C7").Value = 0.04041408728924978
$ws.Range("D7").Value = 0.02094599159831034
$ws.Range("E7").Value = 0.08941456131567094
$ws.Range("F7").Value = 0.6573145730180414
$ws.Range("G7").Value = 0.002420135187214192
$ws.Range("K7").Value = 0.9516407883212707
$ws.Range("M7").Value = 0.3181338523226458
$ws.Range("O7").Value = 2.244972116523286
$ws.Range("C8").Value = 0.04743923049579735
$ws.Range("D8").Value = 0.02279692617150886
$ws.Range("E8").Value = 0.09433170858888928
$ws.Range("F8").Value = 0.6572675424430159
$ws.Range("G8").Value = 0.002415919638966352
$ws.Range("K8").Value = 1.140072074459738
$ws.Range("M8").Value = 0.3705927582217186
$ws.Range("O8").Value = 2.22446877231863
$ws.Range("C9").Value = 0.06133494551016838
$ws.Range("D9").Value = 0.02640561363645588
$ws.Range("E9").Value = 0.104904552020173
$ws.Range("F9").Value = 0.6622212676352106
$ws.Range("G9").Value = 0.002408471113721576
$ws.Range("K9").Value = 1.507573514231581
$ws.Range("M9").Value = 0.4739191239616787
$ws.Range("O9").Value = 2.201564302808521
$ws.Range("C10").Value = 0.07162077706347247
$ws.Range("D10").Value = 0.02904294357063719
$ws.Range("E10").Value = 0.1132614736612041
$ws.Range("F10").Value = 0.6689749086682752
$ws.Range("G10").Value = 0.002403492224868034
$ws.Range("K10").Value = 1.776296137607062
$ws.Range("M10").Value = 0.5501249862443984
$ws.Range("O10").Value = 2.195413457230671
$ws.Range("C11").Value = 0.07631759914855252
$ws.Range("D11").Value = 0.03023947668888383
$ws.Range("E11").Value = 0.1171959390066633
$ws.Range("F11").Value = 0.6727317843004386
$ws.Range("G11").Value = 0.00240133325348157
$ws.Range("K11").Value = 1.898257636291532
$ws.Range("M11").Value = 0.5848602013224138
$ws.Range("O11").Value = 2.194961624570595
$ws.Range("C12").Value = 0.07809876378640013
$ws.Range("D12").Value = 0.03069208872567231
$ws.Range("E12").Value = 0.118705294897488
$ws.Range("F12").Value = 0.6742534969284293
$ws.Range("G12").Value = 0.002400530858275733
$ws.Range("K12").Value = 1.944399426259793
$ws.Range("M12").Value = 0.5980235567196388
$ws.Range("O12").Value = 2.195130126496537
$ws.Range("C13").Value = 0.07771504291413578
$ws.Range("D13").Value = 0.03059463282788499
$ws.Range("E13").Value = 0.1183793569884628
$ws.Range("F13").Value = 0.6739213530747321
$ws.Range("G13").Value = 0.002400702995453492
$ws.Range("K13").Value = 1.934463880525982
$ws.Range("M13").Value = 0.5951881531015033
$ws.Range("O13").Value = 2.19507869742344
$ws.Range("C14").Value = 0.07646408476485078
$ws.Range("D14").Value = 0.03027672331810294
$ws.Range("E14").Value = 0.1173197225972658
$ws.Range("F14").Value = 0.6728549875524976
$ws.Range("G14").Value = 0.002401266936455001
$ws.Range("K14").Value = 1.902054609096638
$ws.Range("M14").Value = 0.5859429604373787
$ws.Range("O14").Value = 2.194968668510484
$ws.Range("C15").Value = 0.0756981736152369
$ws.Range("D15").Value = 0.03008193002827397
$ws.Range("E15").Value = 0.1166732106444925
$ws.Range("F15").Value = 0.672214727479485
$ws.Range("G15").Value = 0.002401614338927327
$ws.Range("K15").Value = 1.882197420299292
$ws.Range("M15").Value = 0.5802812999633176
$ws.Range("O15").Value = 2.194945564534976
$ws.Range("C16").Value = 0.07131418819608371
$ws.Range("D16").Value = 0.0289646804267889
$ws.Range("E16").Value = 0.1130070508168757
$ws.Range("F16").Value = 0.6687432170705705
$ws.Range("G16").Value = 0.002403635444297589
$ws.Range("K16").Value = 1.768319836772775
$ws.Range("M16").Value = 0.547856341548524
$ws.Range("O16").Value = 2.195490387058669
$ws.Range("C17").Value = 0.06862932031492619
$ws.Range("D17").Value = 0.02827844337577545
$ws.Range("E17").Value = 0.1107922605967318
$ws.Range("F17").Value = 0.6667893438484782
$ws.Range("G17").Value = 0.002404902410028477
$ws.Range("K17").Value = 1.698386122592979
$ws.Range("M17").Value = 0.52798228955632
$ws.Range("O17").Value = 2.196427118945252
$ws.Range("C18").Value = 0.06708672394314874
$ws.Range("D18").Value = 0.02788343831117146
$ws.Range("E18").Value = 0.1095308559679538
$ws.Range("F18").Value = 0.665729929447366
$ws.Range("G18").Value = 0.00240564111204667
$ws.Range("K18").Value = 1.658135687675326
$ws.Range("M18").Value = 0.5165577350129098
$ws.Range("O18").Value = 2.19718660689378
$ws.Range("C19").Value = 0.06656471304481215
$ws.Range("D19").Value = 0.02774964580201811
$ws.Range("E19").Value = 0.109105899428009
$ws.Range("F19").Value = 0.6653822721307279
$ws.Range("G19").Value = 0.002405892939762566
$ws.Range("K19").Value = 1.644503103032037
$ws.Range("M19").Value = 0.5126906895090713
$ws.Range("O19").Value = 2.197481595112549
$ws.Range("C20").Value = 0.06891495618550891
$ws.Range("D20").Value = 0.02835152570465738
$ws.Range("E20").Value = 0.1110267341768889
$ws.Range("F20").Value = 0.6669906671535983
$ws.Range("G20").Value = 0.002404766507261728
$ws.Range("K20").Value = 1.705833433032012
$ws.Range("M20").Value = 0.5300972463160463
$ws.Range("O20").Value = 2.196304544065356
$ws.Range("C21").Value = 0.07683145124927648
$ws.Range("D21").Value = 0.03037011454718908
$ws.Range("E21").Value = 0.1176304319548436
$ws.Range("F21").Value = 0.6731655114183184
$ws.Range("G21").Value = 0.002401100882642782
$ws.Range("K21").Value = 1.91157516267765
$ws.Range("M21").Value = 0.5886582294282334
$ws.Range("O21").Value = 2.194991752317065
$ws.Range("C22").Value = 0.0820203900707952
$ws.Range("D22").Value = 0.03168651537819756
$ws.Range("E22").Value = 0.1220599191367313
$ws.Range("F22").Value = 0.6777788447608799
$ws.Range("G22").Value = 0.002398793517183656
$ws.Range("K22").Value = 2.045791395819549
$ws.Range("M22").Value = 0.626988826016742
$ws.Range("O22").Value = 2.19611425154099
$ws.Range("C23").Value = 0.07924956918910198
$ws.Range("D23").Value = 0.03098419883497172
$ws.Range("E23").Value = 0.1196853076249909
$ws.Range("F23").Value = 0.6752635564856462
$ws.Range("G23").Value = 0.00240001694372041
$ws.Range("K23").Value = 1.974180954878705
$ws.Range("M23").Value = 0.606525794804412
$ws.Range("O23").Value = 2.195333205381388
$ws.Range("C24").Value = 0.06878581719088572
$ws.Range("D24").Value = 0.02831848667175052
$ws.Range("E24").Value = 0.1109206916087473
$ws.Range("F24").Value = 0.6668994499066514
$ws.Range("G24").Value = 0.002404827916670972
$ws.Range("K24").Value = 1.702466642811885
$ws.Range("M24").Value = 0.5291410703148642
$ws.Range("O24").Value = 2.196359272077189
$ws.Range("C25").Value = 0.05756260375602551
$ws.Range("D25").Value = 0.0254317391691643
$ws.Range("E25").Value = 0.1019424608607693
$ws.Range("F25").Value = 0.6603370198988614
$ws.Range("G25").Value = 0.002410399090295086
$ws.Range("K25").Value = 1.408375745699914
$ws.Range("M25").Value = 0.4459164284664467
$ws.Range("O25").Value = 2.205894878687246
